$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Vendas" ---
$ws1 = $wb.Worksheets.Item(1)

# Header update
$ws1.Range("F1").Value = "Data Venda"

# Row 2 (existing) updated with new data
$ws1.Range("A2").Value = "devic"
$ws1.Range("B2").Value = "Cimento nassau"
$ws1.Range("C2").Value = 925.0
$ws1.Range("D2").Value = 25
$ws1.Range("E2").Value = "pix"
$ws1.Range("F2").Value = "21/03/2023"

# Row 3 (existing) updated with new data
$ws1.Range("A3").Value = "devic"
$ws1.Range("B3").Value = "Cimento nassau"
$ws1.Range("C3").Value = 23.0
$ws1.Range("D3").Value = 1
$ws1.Range("E3").Value = "pix"
$ws1.Range("F3").Value = "23/01/2023"

# Row 4 (new)
$ws1.Range("A4").Value = "devic"
$ws1.Range("B4").Value = "Cimento nassau"
$ws1.Range("C4").Value = 70.0
$ws1.Range("D4").Value = 2
$ws1.Range("E4").Value = "pix"
$ws1.Range("F4").Value = "28/03/2023"

# Row 5 (new)
$ws1.Range("A5").Value = "devic"
$ws1.Range("B5").Value = "Cimento nassau"
$ws1.Range("C5").Value = 35.0
$ws1.Range("D5").Value = 1
$ws1.Range("E5").Value = "pix"
$ws1.Range("F5").Value = "29/04/2023"

# --- Sheet 2: "Ganhos" ---
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A2").Value = 35.0
$ws2.Range("B2").Value = 1053.0
$ws2.Range("C2").Value = 1053.0
